# Apply cryptos-list price/volume refresh (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cell (D4 / TetherUSD) is untouched by this update and keeps the
# sheet default (un-styled) cell format -- reuse its .Style after forcing
# a text NumberFormat so numeric-looking strings (e.g. "1.00", "4.30") are
# written as literal text (matching the workbook's inlineStr cells) instead
# of being auto-coerced into numbers, while leaving the cell style untouched.
$plainStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = "59.542.50"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "2.513.08"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.77"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("E6").Value = "  -4.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.565"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("D9").Value = "2.514.59"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.39"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  -3.76%  "
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").Value = "2.963.78"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.38"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "59.402.43"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "2.511.35"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("E19").Value = "  -1.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.30"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.66"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.87"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.53"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.423"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  -3.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.171"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = "  +3.52%  "
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.81"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").Value = "0.0₃0779"
$ws.Range("E30").Value = "  -2.82%  "
$ws.Range("E31").Value = "  -1.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.43"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("E34").Value = "  -3.83%  "
$ws.Range("E35").Value = "  -7.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.50"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.23"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  -4.65%  "
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.94"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.813"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  -4.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.24"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  -8.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "281.29"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  -5.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.87"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.597"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.24"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("E50").Value = "  -2.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.87"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  -2.99%  "
